{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Paragraph 1 (index 1): \"Generated: {{DATE}}\" -> \"Date | Activity\"\nparagraphs.items[1].clear();\nparagraphs.items[1].insertText(\"Date | Activity\", Word.InsertLocation.start);\n\n// Paragraph 2 (index 2): empty paragraph -> \"----------------------------\"\nparagraphs.items[2].clear();\nparagraphs.items[2].insertText(\"----------------------------\", Word.InsertLocation.start);\n\n// Paragraph 3 (index 3): \"{{SUMMARY}}\" -> \"{{TABLE}}\"\nparagraphs.items[3].clear();\nparagraphs.items[3].insertText(\"{{TABLE}}\", Word.InsertLocation.start);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$p2 = $d.Paragraphs(2).Range\n$p2.Text = \"Date | Activity\"\n\n$p3 = $d.Paragraphs(3).Range\n$p3.Text = \"----------------------------\"\n\n$p4 = $d.Paragraphs(4).Range\n$p4.Text = \"{{TABLE}}\"\n"}
